$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new data row (19) for period "2509", duplicating the last
#    existing data row (18, period "2508" - Darwis David Martinez Barboza)
#    and shifting the footer rows (old 23/24) down to 24/25.
# ---------------------------------------------------------------------------
$ws.Rows.Item(18).Copy() | Out-Null
$ws.Rows.Item(19).Insert() | Out-Null
$excel.CutCopyMode = $false

# New row 19 now has the same contents/format as old row 18 (period 2508).
# Update its period value to 2509 (new unique period added to shared strings).
$ws.Range("E19").Value = "2509"

# The copy/insert above did not carry the thin grid borders onto the new
# row, so reapply them explicitly (all sides + inside vertical) to match
# the rest of the data table (B:J).
$rng19 = $ws.Range("B19:J19")
$rng19.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$rng19.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$rng19.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
$rng19.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$rng19.Borders.Item(11).LineStyle = 1  # xlInsideVertical

# Row 18 (now a "middle" row instead of the table's last row) keeps its
# thin borders already (no visible change needed there).

# Center-align the "Periodo Mora" column for all data rows, including the
# newly added one.
$ws.Range("E16:E19").HorizontalAlignment = -4108   # xlCenter

# ---------------------------------------------------------------------------
# 2) Update the summary figures: one more period now on record (3 -> 4) and
#    the total overdue amount grows by the new period's "Valor Mora" (56940).
# ---------------------------------------------------------------------------
$ws.Range("F13").Value = 4
$ws.Range("E11").Value = 175660

Write-Output "edit applied"
